$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3300, 3583, 3785, 4240, 4515, 4515, 4671, 4671, 4998, 5031, 5088, 5088, 5088, 5185)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
